$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Rv2174", 3, "mptA Rv2174", "FUNCTION: Involved in the latter stages of the biosynthesis of the alpha-(1->6) mannan core of lipomannan (LM). Catalyzes the addition of alpha-(1->6)-mannose residue. {ECO:0000269|PubMed:17714444}.", 25),
    @("Rv1401", 2, "Rv1401 MTCY21B4.18", "", 25),
    @("Rv1045", 1, "Rv1045", "", 25),
    @("Rv0266c", 1, "oplA Rv0266c", "", 25),
    @("Rv3237c", 1, "Rv3237c", "", 25),
    @("Rv3446c", 1, "Rv3446c", "", 25),
    @("Rv3050c", 1, "Rv3050c", "", 25),
    @("Rv2853", 1, "PE_PGRS48 Rv2853", "", 25),
    @("Rv2943", 1, "Rv2943", "", 25),
    @("Rv2943", 1, "Rv2943A", "", 25),
    @("Rv0030", 1, "Rv0030 MTCY10H4.30", "", 25),
    @("Rv3626c", 1, "Rv3626c", "", 25),
    @("Rv0885", 1, "Rv0885 MTCY31.13", "", 25),
    @("Rv3747", 1, "Rv3747", "", 25)
)

$rowIndex = 6
foreach ($rec in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $rec[0]
    $ws.Cells.Item($rowIndex, 2).Value = $rec[1]
    $ws.Cells.Item($rowIndex, 3).Value = $rec[2]
    if ($rec[3] -ne "") {
        $ws.Cells.Item($rowIndex, 4).Value = $rec[3]
    }
    $ws.Cells.Item($rowIndex, 5).Value = $rec[4]
    $rowIndex++
}
